$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value for semana 5 (row 6)
$ws.Range("B6").Value = 499

# Add new row for semana 6 (row 7)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 2
